# LoginData.xlsx edit:
#   - rename Sheet3 -> customer, Sheet4 -> vendor
#   - fix typo in the customer sheet's bankIfscCode value (IEKE0/92812 -> IEKE092812)
#   - move the selected cell on the customer sheet from N5 to N7
#   - make "vendor" the active/selected sheet (was "customer")

$wb = $excel.ActiveWorkbook

# Sheets 3 and 4 are named "Sheet3"/"Sheet4"; rename them to their real names.
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)
$ws3.Name = "customer"
$ws4.Name = "vendor"

# Correct the stray "/" in the IFSC code stored in N3 on the "customer" sheet.
$ws3.Range("N3").Value = "IEKE092812"

# The "customer" sheet was the active tab with N5 selected; update the selection
# to N7, then make "vendor" the active tab (it becomes tabSelected="true",
# "customer" becomes tabSelected="false").
[void]$ws3.Activate()
[void]$ws3.Range("N7").Select()

[void]$ws4.Activate()
